# Natmi following Dr Hou advice
# Update Slitrk1-Ptprs LR-pair results: recompute metrics for FAPs/sCs
# sending clusters and add a new "ECs" sending-cluster block (rows 2-16).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Slitrk1"
$ws.Range("C2").Value = "Ptprs"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.018579
$ws.Range("H2").Value = 0.055737
$ws.Range("I2").Value = 0.1519646868933627
$ws.Range("J2").Value = 0.1519646868933627
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 4.260928666666667
$ws.Range("N2").Value = 12.782786
$ws.Range("O2").Value = 0.05853043679913345
$ws.Range("P2").Value = 0.05853043679913345
$ws.Range("Q2").Value = 0.079163793698
$ws.Range("R2").Value = 0.712474143282
$ws.Range("S2").Value = 0.008894559501912071
$ws.Range("T2").Value = 0.008894559501912069

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Slitrk1"
$ws.Range("C3").Value = "Ptprs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.018579
$ws.Range("H3").Value = 0.055737
$ws.Range("I3").Value = 0.1519646868933627
$ws.Range("J3").Value = 0.1519646868933627
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 47.585289
$ws.Range("N3").Value = 142.755867
$ws.Range("O3").Value = 0.6536574461270807
$ws.Range("P3").Value = 0.6536574461270805
$ws.Range("Q3").Value = 0.884087084331
$ws.Range("R3").Value = 7.956783758979
$ws.Range("S3").Value = 0.09933284913621693
$ws.Range("T3").Value = 0.0993328491362169

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Slitrk1"
$ws.Range("C4").Value = "Ptprs"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.018579
$ws.Range("H4").Value = 0.055737
$ws.Range("I4").Value = 0.1519646868933627
$ws.Range("J4").Value = 0.1519646868933627
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.555122999999999
$ws.Range("N4").Value = 16.665369
$ws.Range("O4").Value = 0.07630819502014176
$ws.Range("P4").Value = 0.07630819502014174
$ws.Range("Q4").Value = 0.103208630217
$ws.Range("R4").Value = 0.928877671953
$ws.Range("S4").Value = 0.0115961509636335
$ws.Range("T4").Value = 0.0115961509636335

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Slitrk1"
$ws.Range("C5").Value = "Ptprs"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.018579
$ws.Range("H5").Value = 0.055737
$ws.Range("I5").Value = 0.1519646868933627
$ws.Range("J5").Value = 0.1519646868933627
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 7.6722
$ws.Range("N5").Value = 23.0166
$ws.Range("O5").Value = 0.1053895177178852
$ws.Range("P5").Value = 0.1053895177178852
$ws.Range("Q5").Value = 0.1425418038
$ws.Range("R5").Value = 1.2828762342
$ws.Range("S5").Value = 0.01601548506184093
$ws.Range("T5").Value = 0.01601548506184092

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Slitrk1"
$ws.Range("C6").Value = "Ptprs"
$ws.Range("D6").Value = "sCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.018579
$ws.Range("H6").Value = 0.055737
$ws.Range("I6").Value = 0.1519646868933627
$ws.Range("J6").Value = 0.1519646868933627
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.724970666666667
$ws.Range("N6").Value = 23.174912
$ws.Range("O6").Value = 0.106114404335759
$ws.Range("P6").Value = 0.106114404335759
$ws.Range("Q6").Value = 0.143522230016
$ws.Range("R6").Value = 1.291700070144
$ws.Range("S6").Value = 0.01612564222975931
$ws.Range("T6").Value = 0.0161256422297593

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Slitrk1"
$ws.Range("C7").Value = "Ptprs"
$ws.Range("D7").Value = "ECs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.08670066666666666
$ws.Range("H7").Value = 0.260102
$ws.Range("I7").Value = 0.7091576329967064
$ws.Range("J7").Value = 0.7091576329967064
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.260928666666667
$ws.Range("N7").Value = 12.782786
$ws.Range("O7").Value = 0.05853043679913345
$ws.Range("P7").Value = 0.05853043679913345
$ws.Range("Q7").Value = 0.3694253560191111
$ws.Range("R7").Value = 3.324828204172
$ws.Range("S7").Value = 0.0415073060187368
$ws.Range("T7").Value = 0.0415073060187368

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Slitrk1"
$ws.Range("C8").Value = "Ptprs"
$ws.Range("D8").Value = "FAPs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.08670066666666666
$ws.Range("H8").Value = 0.260102
$ws.Range("I8").Value = 0.7091576329967064
$ws.Range("J8").Value = 0.7091576329967064
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 47.585289
$ws.Range("N8").Value = 142.755867
$ws.Range("O8").Value = 0.6536574461270807
$ws.Range("P8").Value = 0.6536574461270805
$ws.Range("Q8").Value = 4.125676279825999
$ws.Range("R8").Value = 37.131086518434
$ws.Range("S8").Value = 0.4635461672861527
$ws.Range("T8").Value = 0.4635461672861526

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Slitrk1"
$ws.Range("C9").Value = "Ptprs"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.08670066666666666
$ws.Range("H9").Value = 0.260102
$ws.Range("I9").Value = 0.7091576329967064
$ws.Range("J9").Value = 0.7091576329967064
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.555122999999999
$ws.Range("N9").Value = 16.665369
$ws.Range("O9").Value = 0.07630819502014176
$ws.Range("P9").Value = 0.07630819502014174
$ws.Range("Q9").Value = 0.4816328675153332
$ws.Range("R9").Value = 4.334695807638
$ws.Range("S9").Value = 0.05411453895873479
$ws.Range("T9").Value = 0.05411453895873478

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Slitrk1"
$ws.Range("C10").Value = "Ptprs"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.08670066666666666
$ws.Range("H10").Value = 0.260102
$ws.Range("I10").Value = 0.7091576329967064
$ws.Range("J10").Value = 0.7091576329967064
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 7.6722
$ws.Range("N10").Value = 23.0166
$ws.Range("O10").Value = 0.1053895177178852
$ws.Range("P10").Value = 0.1053895177178852
$ws.Range("Q10").Value = 0.6651848548
$ws.Range("R10").Value = 5.9866636932
$ws.Range("S10").Value = 0.07473778092747993
$ws.Range("T10").Value = 0.07473778092747992

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Slitrk1"
$ws.Range("C11").Value = "Ptprs"
$ws.Range("D11").Value = "sCs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.08670066666666666
$ws.Range("H11").Value = 0.260102
$ws.Range("I11").Value = 0.7091576329967064
$ws.Range("J11").Value = 0.7091576329967064
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 7.724970666666667
$ws.Range("N11").Value = 23.174912
$ws.Range("O11").Value = 0.106114404335759
$ws.Range("P11").Value = 0.106114404335759
$ws.Range("Q11").Value = 0.6697601067804444
$ws.Range("R11").Value = 6.027840961023999
$ws.Range("S11").Value = 0.07525183980560229
$ws.Range("T11").Value = 0.07525183980560228

# Row 12
$ws.Range("A12").Value = "sCs"
$ws.Range("B12").Value = "Slitrk1"
$ws.Range("C12").Value = "Ptprs"
$ws.Range("D12").Value = "ECs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.016979
$ws.Range("H12").Value = 0.050937
$ws.Range("I12").Value = 0.1388776801099308
$ws.Range("J12").Value = 0.1388776801099308
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 4.260928666666667
$ws.Range("N12").Value = 12.782786
$ws.Range("O12").Value = 0.05853043679913345
$ws.Range("P12").Value = 0.05853043679913345
$ws.Range("Q12").Value = 0.07234630783133332
$ws.Range("R12").Value = 0.6511167704819999
$ws.Range("S12").Value = 0.00812857127848458
$ws.Range("T12").Value = 0.008128571278484579

# Row 13
$ws.Range("A13").Value = "sCs"
$ws.Range("B13").Value = "Slitrk1"
$ws.Range("C13").Value = "Ptprs"
$ws.Range("D13").Value = "FAPs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.016979
$ws.Range("H13").Value = 0.050937
$ws.Range("I13").Value = 0.1388776801099308
$ws.Range("J13").Value = 0.1388776801099308
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 47.585289
$ws.Range("N13").Value = 142.755867
$ws.Range("O13").Value = 0.6536574461270807
$ws.Range("P13").Value = 0.6536574461270805
$ws.Range("Q13").Value = 0.8079506219309999
$ws.Range("R13").Value = 7.271555597378999
$ws.Range("S13").Value = 0.09077842970471106
$ws.Range("T13").Value = 0.09077842970471105

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Slitrk1"
$ws.Range("C14").Value = "Ptprs"
$ws.Range("D14").Value = "M1"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.016979
$ws.Range("H14").Value = 0.050937
$ws.Range("I14").Value = 0.1388776801099308
$ws.Range("J14").Value = 0.1388776801099308
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.555122999999999
$ws.Range("N14").Value = 16.665369
$ws.Range("O14").Value = 0.07630819502014176
$ws.Range("P14").Value = 0.07630819502014174
$ws.Range("Q14").Value = 0.09432043341699997
$ws.Range("R14").Value = 0.8488839007529998
$ws.Range("S14").Value = 0.01059750509777346
$ws.Range("T14").Value = 0.01059750509777346

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Slitrk1"
$ws.Range("C15").Value = "Ptprs"
$ws.Range("D15").Value = "M2"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.016979
$ws.Range("H15").Value = 0.050937
$ws.Range("I15").Value = 0.1388776801099308
$ws.Range("J15").Value = 0.1388776801099308
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 7.6722
$ws.Range("N15").Value = 23.0166
$ws.Range("O15").Value = 0.1053895177178852
$ws.Range("P15").Value = 0.1053895177178852
$ws.Range("Q15").Value = 0.1302662838
$ws.Range("R15").Value = 1.1723965542
$ws.Range("S15").Value = 0.01463625172856435
$ws.Range("T15").Value = 0.01463625172856435

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Slitrk1"
$ws.Range("C16").Value = "Ptprs"
$ws.Range("D16").Value = "sCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.016979
$ws.Range("H16").Value = 0.050937
$ws.Range("I16").Value = 0.1388776801099308
$ws.Range("J16").Value = 0.1388776801099308
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 7.724970666666667
$ws.Range("N16").Value = 23.174912
$ws.Range("O16").Value = 0.106114404335759
$ws.Range("P16").Value = 0.106114404335759
$ws.Range("Q16").Value = 0.1311622769493333
$ws.Range("R16").Value = 1.180460492544
$ws.Range("S16").Value = 0.0147369223003974
$ws.Range("T16").Value = 0.0147369223003974

